# task 4, fixed normal distribution
# Updates the "True (Q_star)" (row 3) and "Absolute Error" (row 4) rows
# on each of the 5 result sheets with corrected values.

$wb = $excel.ActiveWorkbook

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

$data = @{
    "Normal (Correct)" = @{
        3 = @(23.6125159492159, 24.22894706844462, 24.35516669106373, 24.79639007135757, 25.01058790125472, 25.31578976705599, 25.60961652558548, 25.97094883707175, 26.06482742811744)
        4 = @(0.2243101137636785, 0.05137388192035885, 0.004057526163968816, 0.05859032771158823, 0.01058790125471987, 0.05358951070197193, 0.03115925718682178, 0.1485220235960156, 0.09834650890298491)
    }
    "Poisson n=10" = @{
        3 = @(21, 20, 19, 22, 25, 26, 25, 30, 33)
        4 = @(2.836826062979579, 4.177573186524263, 5.359224217227698, 2.737799743645979, 0, 0.7377997436459793, 0.6407757827723017, 4.177573186524263, 6.836826062979579)
    }
    "Poisson n=200" = @{
        3 = @(16, 18, 18, 22, 24, 28, 32, 33, 38)
        4 = @(7.836826062979579, 6.177573186524263, 6.359224217227698, 2.737799743645979, 1, 2.737799743645979, 6.359224217227698, 7.177573186524263, 11.83682606297958)
    }
    "Exponential n=10" = @{
        3 = @(0.3530271770614896, 0.08317540839287309, 1.895093424053536, 5.182297750219801, 20.44251913370648, 36.2338528758555, 45.38043345149831, 59.04172972873226, 66.69359476997204)
        4 = @(23.48379888591809, 24.09439777813139, 22.46413079317416, 19.55550199342618, 4.557480866293524, 10.97165261950148, 19.739657668726, 33.21930291525653, 40.53042083295162)
    }
    "Exponential n=200" = @{
        3 = @(0.1095801702055237, 1.145873609188977, 2.823217022293488, 8.979461359037044, 20.43130362617683, 23.33631828311393, 62.41232333411797, 70.55005985866988, 107.3466431768455)
        4 = @(23.72724589277406, 23.03169957733529, 21.53600719493421, 15.75833838460894, 4.568696373823169, 1.925881973240095, 36.77154755134567, 44.72763304519415, 81.18346923982506)
    }
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $data[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $values = $rows[$rowNum]
        for ($i = 0; $i -lt $cols.Length; $i++) {
            $cellRef = "$($cols[$i])$rowNum"
            $ws.Range($cellRef).Value = $values[$i]
        }
    }
}
